# Daily-scrape update: two existing listings replaced (Tanta / Monkey Brew
# postings superseded by Shebeen El-Kom / European Hospital postings) and a
# brand-new third listing appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    <#
        Writes a value into a cell while forcing Excel's COM "smart" type
        inference to treat it as plain text (needed for the numeric-looking
        opportunity IDs in column A) and then resets the cell style back to
        "Normal" so no stray number-format style is left attached to the
        cell - it only ever carries the literal string, same as before.
    #>
    param($range, $value)

    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---- Row 2: opportunity 1327758 -> 1327807 -----------------------------
Set-TextValue $ws.Range("A2") "1327807"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327807"
$ws.Range("C2").Value = "Marketing Specialist"
$ws.Range("D2").Value = "Shebeen El-Kom, Qism Shebeen El-Kom, Shibin el Kom, Menofia Governorate, Egypt"
$ws.Range("F2").Value = "1 applicant"
$ws.Range("H2").Value = "European Hospital"

# ---- Row 3: opportunity 1327757 -> 1327806 -----------------------------
Set-TextValue $ws.Range("A3") "1327806"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327806"
$ws.Range("C3").Value = "Office manager"
$ws.Range("D3").Value = "Shebeen El-Kom, Qism Shebeen El-Kom, Shibin el Kom, Menofia Governorate, Egypt"
$ws.Range("H3").Value = "European Hospital"

# ---- Row 4: brand-new opportunity 1326583 ------------------------------
Set-TextValue $ws.Range("A4") "1326583"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326583"
$ws.Range("C4").Value = "Software Developer"
$ws.Range("D4").Value = "2750 Cascais, Portugal"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "107 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "MQ Capital"

# ---- Column width tweaks (C, D, F, H got wider) ------------------------
# ColumnWidth is in "characters" and Excel stores the OOXML <col width>
# with the familiar +5/6 padding offset, so back that off here to land on
# the exact stored widths from the diff (23 / 81 / 17 / 20).
$ws.Columns.Item(3).ColumnWidth = 23 - 5/6
$ws.Columns.Item(4).ColumnWidth = 81 - 5/6
$ws.Columns.Item(6).ColumnWidth = 17 - 5/6
$ws.Columns.Item(8).ColumnWidth = 20 - 5/6
